$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The target cells store values that look numeric/percent but are kept as
# plain text strings in the source workbook, so force text format per-cell
# before assigning to avoid Excel auto-converting them to numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "277.92"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "1.30%"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "0.90%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.855"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "1.23%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06407"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "1.97%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.989"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "1.09%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.201"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-8.00%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8782"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "0.92%"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-1.92%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.05168"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "3.32%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07505"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "0.32%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.02931"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "1.69%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.08971"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.92%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001568"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.04%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0006381"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.64%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006081"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "3.51%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.474"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.60%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.306"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-0.16%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.246"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-1.69%"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-2.01%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1323"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "0.42%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.904"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-0.53%"
$ws.Range("B23").NumberFormat = "@"
$ws.Range("B23").Value = "CoinExToken"
$ws.Range("C23").NumberFormat = "@"
$ws.Range("C23").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04410"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "0.78%"
$ws.Range("B24").NumberFormat = "@"
$ws.Range("B24").Value = "ZBToken"
$ws.Range("C24").NumberFormat = "@"
$ws.Range("C24").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.1506"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "9.11%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.001178"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "0.76%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.003896"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "2.23%"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0001182"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "-1.61%"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "1.71%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04079"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "0.14%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006809"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-3.58%"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "0.47%"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-7.00%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01120"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-0.32%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005360"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "3.55%"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "4.79%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-19.52%"
